$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.864.12'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.61%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.391.30'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.46%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '560.49'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.72%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.80'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.61%  '

$ws.Range("E7").Value = '  +2.58%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.382.88'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.61%  '

$ws.Range("E9").Value = '  -0.02%  '

$ws.Range("E10").Value = '  +10.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.632'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.22%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.61'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.46%  '

$ws.Range("E13").Value = '  +4.47%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.14'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.67%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.936.09'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.26%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '18.27'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.20%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.390.95'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.88%  '

$ws.Range("E18").Value = '  +0.50%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.90'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.86%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '64.917.18'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.63%  '

$ws.Range("E21").Value = '  +2.55%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '469.25'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +16.33%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +16.69%  '

$ws.Range("E24").Value = '  +1.39%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.31'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.31%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.65'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.78%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.83'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.48%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.86'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.34%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.84'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.39%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.64'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.08%  '

$ws.Range("E31").Value = '  +4.02%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.55'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.84%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '579.65'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.56%  '

$ws.Range("E34").Value = '  +2.28%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '60.10'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.93%  '

$ws.Range("E36").Value = '  +0.12%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.141'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.77%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '35.93'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.65%  '

$ws.Range("E39").Value = '  +1.87%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0755'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.06%  '

$ws.Range("E41").Value = '  +1.09%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.100.35'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.53%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.24%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.87'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.06%  '

$ws.Range("E45").Value = '  +1.90%  '

$ws.Range("E46").Value = '  +2.34%  '

$ws.Range("E47").Value = '  -0.99%  '

$ws.Range("E48").Value = '  +4.27%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.56'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.16%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '137.81'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.92%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.39'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.95%  '
